# Updated cryptos list (Price/Volume(1h) refresh), mirrors the upstream
# GitHub Actions commit. D-column prices are entered with a leading
# apostrophe so Excel stores them as literal text (matching the source
# file's inlineStr cells) instead of coercing look-alike numbers/dates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''30.140.31'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '''1.915.07'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('D4').Value = '''0.9997'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''0.7910'
$ws.Range('E5').Value = '  +6.58%  '
$ws.Range('D6').Value = '''242.73'
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('D7').Value = '''1.000'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '''0.3183'
$ws.Range('E8').Value = '  +2.95%  '
$ws.Range('D9').Value = '''26.36'
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('D10').Value = '''0.06960'
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('D11').Value = '''0.08005'
$ws.Range('E11').Value = '  -0.89%  '
$ws.Range('D12').Value = '''0.7522'
$ws.Range('E12').Value = '  -2.44%  '
$ws.Range('D13').Value = '''1.915.07'
$ws.Range('E13').Value = '  -0.95%  '
$ws.Range('D14').Value = '''5.227'
$ws.Range('E14').Value = '  -1.86%  '
$ws.Range('D15').Value = '''93.51'
$ws.Range('E15').Value = '  +1.34%  '
$ws.Range('D16').Value = '''30.166.56'
$ws.Range('D17').Value = '''14.07'
$ws.Range('E17').Value = '  -1.46%  '
$ws.Range('D18').Value = '''5.994'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('D19').Value = '''249.01'
$ws.Range('E19').Value = '  +3.66%  '
$ws.Range('D20').Value = '''0.000007828'
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('D21').Value = '''0.9996'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').Value = '''2.152.80'
$ws.Range('E22').Value = '  -1.86%  '
$ws.Range('D23').Value = '''1.001'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = '''6.977'
$ws.Range('E24').Value = '  -1.94%  '
$ws.Range('D25').Value = '''169.21'
$ws.Range('E25').Value = '  +1.44%  '
$ws.Range('D26').Value = '''9.322'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').Value = '''0.1386'
$ws.Range('E27').Value = '  +8.69%  '
$ws.Range('D28').Value = '''18.99'
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('D29').Value = '''2.056'
$ws.Range('E29').Value = '  +0.51%  '
$ws.Range('D30').Value = '''1.391'
$ws.Range('E30').Value = '  +3.00%  '
$ws.Range('D31').Value = '''1.528'
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('D32').Value = '''4.367'
$ws.Range('E32').Value = '  +0.64%  '
$ws.Range('D33').Value = '''4.124'
$ws.Range('E33').Value = '  +0.98%  '
$ws.Range('D34').Value = '''0.05395'
$ws.Range('E34').Value = '  +4.58%  '
$ws.Range('D35').Value = '''1.267'
$ws.Range('E35').Value = '  -3.26%  '
$ws.Range('D36').Value = '''0.7399'
$ws.Range('E36').Value = '  -1.31%  '
$ws.Range('E37').Value = '  +0.30%  '
$ws.Range('D38').Value = '''0.01932'
$ws.Range('E38').Value = '  -1.44%  '
$ws.Range('D39').Value = '''2.790'
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('D40').Value = '''6.185'
$ws.Range('E40').Value = '  -2.71%  '
$ws.Range('D41').Value = '''0.4460'
$ws.Range('E41').Value = '  -0.96%  '
$ws.Range('D42').Value = '''72.74'
$ws.Range('E42').Value = '  -2.40%  '
$ws.Range('D43').Value = '''1.906'
$ws.Range('E43').Value = '  -3.79%  '
$ws.Range('D44').Value = '''1.000'
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('D45').Value = '''0.8346'
$ws.Range('E45').Value = '  -0.56%  '
$ws.Range('D46').Value = '''7.613'
$ws.Range('E46').Value = '  -1.45%  '
$ws.Range('D47').Value = '''9.855'
$ws.Range('E47').Value = '  -1.13%  '
$ws.Range('D48').Value = '''100.59'
$ws.Range('E48').Value = '  -1.19%  '
$ws.Range('D49').Value = '''2.060.72'
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('D50').Value = '''965.46'
$ws.Range('E50').Value = '  +4.22%  '
$ws.Range('D51').Value = '''36.53'
$ws.Range('E51').Value = '  -0.69%  '
